$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O: header year 2021 (row 4) and value (row 5)
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 4.0999999999999996

# L5 value change 1.7 -> 1.6
$ws.Range("L5").Value = 1.6

# N5 value change 1.6 -> 3.1
$ws.Range("N5").Value = 3.1

# Update the selected cell to P4 (was P6)
$ws.Range("P4").Select()
